$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A31").Value = "Answer:"
$ws.Range("B31").Value = '${answerToLifeTheUniverseAndEverything}'

$ws.Range("A32").Value = "Pick A Card:"
$ws.Range("B32").Value = '${jett:pickACard()}'
